$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new problem entry (P1088) as row 31
$ws.Range("A31").Value = "P1088"
$ws.Range("B31").Value = "火星人"
$ws.Range("C31").Value = "AC"
$ws.Range("D31").Value = "普及/提高-"
$ws.Range("E31").Value = "数学"
$ws.Range("F31").Value = "使用STL算法next_permutation可以水过"
$ws.Range("F31").WrapText = $true
$ws.Range("G31").Value = "2019-11-25"
$ws.Range("H31").Value = "2019-11-25"

# Update the view so the newly added row is visible / selected, matching
# the author re-saving after scrolling down to the new entry.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 18
$aw.ScrollColumn = 1
[void]$ws.Range("H32").Select()

Write-Host "Row 31 (P1088) added."
